$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3845463333333334
$ws.Range("H2").Value = 1.153639
$ws.Range("I2").Value = 0.1984850200147207
$ws.Range("J2").Value = 0.1984850200147207
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.841766333333334
$ws.Range("N2").Value = 11.525299
$ws.Range("O2").Value = 0.04788297632726236
$ws.Range("P2").Value = 0.04788297632726236
$ws.Range("Q2").Value = 1.477337157006778
$ws.Range("R2").Value = 13.296034413061
$ws.Range("S2").Value = 0.009504053514681069
$ws.Range("T2").Value = 0.009504053514681067
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3845463333333334
$ws.Range("H3").Value = 1.153639
$ws.Range("I3").Value = 0.1984850200147207
$ws.Range("J3").Value = 0.1984850200147207
$ws.Range("O3").Value = 0.5601341401483774
$ws.Range("P3").Value = 0.5601341401483775
$ws.Range("Q3").Value = 17.281861772617
$ws.Range("R3").Value = 155.536755953553
$ws.Range("S3").Value = 0.1111782360182791
$ws.Range("T3").Value = 0.1111782360182791
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3845463333333334
$ws.Range("H4").Value = 1.153639
$ws.Range("I4").Value = 0.1984850200147207
$ws.Range("J4").Value = 0.1984850200147207
$ws.Range("M4").Value = 31.44972933333333
$ws.Range("N4").Value = 94.349188
$ws.Range("O4").Value = 0.3919828835243602
$ws.Range("P4").Value = 0.3919828835243602
$ws.Range("Q4").Value = 12.09387809945911
$ws.Range("R4").Value = 108.844902895132
$ws.Range("S4").Value = 0.07780273048176058
$ws.Range("T4").Value = 0.07780273048176058
$ws.Range("I5").Value = 0.5733580031870772
$ws.Range("J5").Value = 0.5733580031870772
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.841766333333334
$ws.Range("N5").Value = 11.525299
$ws.Range("O5").Value = 0.04788297632726236
$ws.Range("P5").Value = 0.04788297632726236
$ws.Range("Q5").Value = 4.267541612524
$ws.Range("R5").Value = 38.407874512716
$ws.Range("S5").Value = 0.02745408769365323
$ws.Range("T5").Value = 0.02745408769365323
$ws.Range("I6").Value = 0.5733580031870772
$ws.Range("J6").Value = 0.5733580031870772
$ws.Range("O6").Value = 0.5601341401483774
$ws.Range("P6").Value = 0.5601341401483775
$ws.Range("R6").Value = 449.294580563868
$ws.Range("S6").Value = 0.3211573921123841
$ws.Range("T6").Value = 0.3211573921123841
$ws.Range("I7").Value = 0.5733580031870772
$ws.Range("J7").Value = 0.5733580031870772
$ws.Range("M7").Value = 31.44972933333333
$ws.Range("N7").Value = 94.349188
$ws.Range("O7").Value = 0.3919828835243602
$ws.Range("P7").Value = 0.3919828835243602
$ws.Range("Q7").Value = 34.935239935888
$ws.Range("R7").Value = 314.417159422992
$ws.Range("S7").Value = 0.2247465233810398
$ws.Range("T7").Value = 0.2247465233810398
$ws.Range("G8").Value = 0.442033
$ws.Range("H8").Value = 1.326099
$ws.Range("I8").Value = 0.2281569767982021
$ws.Range("J8").Value = 0.2281569767982021
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.841766333333334
$ws.Range("N8").Value = 11.525299
$ws.Range("O8").Value = 0.04788297632726236
$ws.Range("P8").Value = 0.04788297632726236
$ws.Range("Q8").Value = 1.698187497622333
$ws.Range("R8").Value = 15.283687478601
$ws.Range("S8").Value = 0.01092483511892806
$ws.Range("T8").Value = 0.01092483511892806
$ws.Range("G9").Value = 0.442033
$ws.Range("H9").Value = 1.326099
$ws.Range("I9").Value = 0.2281569767982021
$ws.Range("J9").Value = 0.2281569767982021
$ws.Range("O9").Value = 0.5601341401483774
$ws.Range("P9").Value = 0.5601341401483775
$ws.Range("Q9").Value = 19.865364827997
$ws.Range("R9").Value = 178.788283451973
$ws.Range("S9").Value = 0.1277985120177142
$ws.Range("T9").Value = 0.1277985120177142
$ws.Range("G10").Value = 0.442033
$ws.Range("H10").Value = 1.326099
$ws.Range("I10").Value = 0.2281569767982021
$ws.Range("J10").Value = 0.2281569767982021
$ws.Range("M10").Value = 31.44972933333333
$ws.Range("N10").Value = 94.349188
$ws.Range("O10").Value = 0.3919828835243602
$ws.Range("P10").Value = 0.3919828835243602
$ws.Range("Q10").Value = 13.90181820640133
$ws.Range("R10").Value = 125.116363857612
$ws.Range("S10").Value = 0.08943362966155981
$ws.Range("T10").Value = 0.08943362966155981
